$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 117
$ws.Range("I11").Value = 117
$ws.Range("K11").Value = 117
$ws.Range("M11").Value = 23

$ws.Range("H42").Value = 22.8
$ws.Range("I42").Value = 22.8
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 68.40000000000001
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 161.6
$ws.Range("N42").ClearContents()

$ws.Range("H55").Value = 366.22223
$ws.Range("J55").Value = 499.5
$ws.Range("L55").Value = 499.5
$ws.Range("N55").Value = -927.5

$ws.Range("H62").Value = 7881.769
$ws.Range("I62").Value = 7050.5
$ws.Range("K62").Value = 7050.5
$ws.Range("M62").Value = -6426.5

$ws.Range("H65").Value = 7881.769
$ws.Range("I65").Value = 7050.5
$ws.Range("K65").Value = 35252.5
$ws.Range("M65").Value = -32132.5

$ws.Range("H103").Value = 607.5714
$ws.Range("I103").Value = 646.4
$ws.Range("J103").Value = 510.5
$ws.Range("K103").Value = 1939.2
$ws.Range("L103").Value = 1531.5
$ws.Range("M103").Value = -1353.2
$ws.Range("N103").Value = -2703.5

$ws.Range("H125").Value = 6371.5
$ws.Range("I125").Value = 7348.2
$ws.Range("J125").Value = 1488
$ws.Range("K125").Value = 66133.8
$ws.Range("L125").Value = 13392
$ws.Range("M125").Value = -63673.8
$ws.Range("N125").Value = -18312

$ws.Range("H137").Value = 11236.912
$ws.Range("I137").Value = 1729.2667
$ws.Range("K137").Value = 5187.800099999999
$ws.Range("M137").Value = -2637.800099999999

$ws.Range("H138").Value = 4834.7856
$ws.Range("J138").Value = 2635.7778
$ws.Range("L138").Value = 7907.3334
$ws.Range("N138").Value = -18187.3334


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15572.947
$ws.Range("I32").Value = 12492.571
$ws.Range("J32").Value = 24198
$ws.Range("K32").Value = 12492.571
$ws.Range("L32").Value = 24198
$ws.Range("M32").Value = -12205.571
$ws.Range("N32").Value = -24772

$ws.Range("H61").Value = 17587.36
$ws.Range("I61").Value = 5567.5
$ws.Range("K61").Value = 5567.5
$ws.Range("M61").Value = -5355.5

$ws.Range("H63").Value = 2565.5
$ws.Range("I63").Value = 2378.6
$ws.Range("K63").Value = 2378.6
$ws.Range("M63").Value = -1692.6

$ws.Range("H66").Value = 2565.5
$ws.Range("I66").Value = 2378.6
$ws.Range("K66").Value = 11893
$ws.Range("M66").Value = -8461

$ws.Range("H122").Value = 4884.5713
$ws.Range("I122").Value = 2548
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 7644
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -5194
$ws.Range("N122").Value = -28900

$ws.Range("H132").Value = 3239307.5
$ws.Range("J132").Value = 12540133
$ws.Range("L132").Value = 37620399
$ws.Range("N132").Value = -37625459

$ws.Range("H136").Value = 17587.36
$ws.Range("I136").Value = 5567.5
$ws.Range("K136").Value = 16702.5
$ws.Range("M136").Value = -14152.5


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6776.533
$ws.Range("I99").Value = 3200.7273
$ws.Range("J99").Value = 8846.736999999999
$ws.Range("K99").Value = 3200.7273
$ws.Range("L99").Value = 8846.736999999999
$ws.Range("M99").Value = -1702.7273
$ws.Range("N99").Value = -11842.737

$ws.Range("H126").Value = 6776.533
$ws.Range("I126").Value = 3200.7273
$ws.Range("J126").Value = 8846.736999999999
$ws.Range("K126").Value = 9602.1819
$ws.Range("L126").Value = 26540.211
$ws.Range("M126").Value = -7132.1819
$ws.Range("N126").Value = -31480.211

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 13211.833
$ws.Range("I132").Value = 3789.7144
$ws.Range("J132").Value = 26402.8
$ws.Range("K132").Value = 11369.1432
$ws.Range("L132").Value = 79208.39999999999
$ws.Range("M132").Value = -8839.143199999999
$ws.Range("N132").Value = -84268.39999999999

$ws.Range("H134").Value = 41675360
$ws.Range("I134").Value = 1895.3636
$ws.Range("K134").Value = 5686.0908
$ws.Range("M134").Value = -3151.0908

$ws.Range("H33").Value = 588.375
$ws.Range("I33").Value = 178.25
$ws.Range("J33").Value = 998.5
$ws.Range("K33").Value = 1069.5
$ws.Range("L33").Value = 5991
$ws.Range("M33").Value = -786.5
$ws.Range("N33").Value = -6557

$ws.Range("H99").Value = 2571.4285
$ws.Range("I99").Value = 1600
$ws.Range("K99").Value = 4800
$ws.Range("M99").Value = -2554

$ws.Range("H114").Value = 840.125
$ws.Range("I114").Value = 933.6
$ws.Range("J114").Value = 684.3333
$ws.Range("K114").Value = 2800.8
$ws.Range("L114").Value = 2052.9999
$ws.Range("M114").Value = 453.1999999999998
$ws.Range("N114").Value = -8560.999899999999

$ws.Range("H131").Value = 1497.8
$ws.Range("I131").Value = 1400
$ws.Range("J131").Value = 1498.7878
$ws.Range("K131").Value = 4200
$ws.Range("L131").Value = 4496.3634
$ws.Range("M131").Value = 840
$ws.Range("N131").Value = -14576.3634


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10598.75
$ws.Range("I20").Value = 6197.5
$ws.Range("K20").Value = 6197.5
$ws.Range("M20").Value = -5952.5

$ws.Range("H24").Value = 2012875.9
$ws.Range("I24").Value = 10004380
$ws.Range("K24").Value = 10004380
$ws.Range("M24").Value = -10004207


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2926.0952
$ws.Range("I22").Value = 2194.7273
$ws.Range("J22").Value = 3730.6
$ws.Range("K22").Value = 2194.7273
$ws.Range("L22").Value = 3730.6
$ws.Range("M22").Value = -1899.7273
$ws.Range("N22").Value = -4320.6

$ws.Range("H27").Value = 2926.0952
$ws.Range("I27").Value = 2194.7273
$ws.Range("J27").Value = 3730.6
$ws.Range("K27").Value = 2194.7273
$ws.Range("L27").Value = 3730.6
$ws.Range("M27").Value = -2087.7273
$ws.Range("N27").Value = -3944.6

$ws.Range("H40").Value = 6175.6113
$ws.Range("I40").Value = 3710
$ws.Range("K40").Value = 3710
$ws.Range("M40").Value = -3574

$ws.Range("H46").Value = 2936.75
$ws.Range("I46").Value = 1649.4
$ws.Range("K46").Value = 1649.4
$ws.Range("M46").Value = -1461.4

$ws.Range("H122").Value = 7791.6665
$ws.Range("I122").Value = 8000
$ws.Range("J122").Value = 7583.3335
$ws.Range("K122").Value = 24000
$ws.Range("L122").Value = 22750.0005
$ws.Range("M122").Value = -21550
$ws.Range("N122").Value = -27650.0005

$ws.Range("H136").Value = 13932.872
$ws.Range("I136").Value = 19582.5
$ws.Range("J136").Value = 11421.926
$ws.Range("K136").Value = 58747.5
$ws.Range("L136").Value = 34265.778
$ws.Range("M136").Value = -56197.5
$ws.Range("N136").Value = -39365.778


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 12750
$ws.Range("J3").Value = 12750
$ws.Range("L3").Value = 12750
$ws.Range("N3").Value = -12978

$ws.Range("H31").Value = 20000
$ws.Range("J31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("N31").Value = -20696

$ws.Range("H62").Value = 3387.1428
$ws.Range("I62").Value = 2986.6667
$ws.Range("J62").Value = 3687.5
$ws.Range("K62").Value = 2986.6667
$ws.Range("L62").Value = 3687.5
$ws.Range("M62").Value = -2362.6667
$ws.Range("N62").Value = -4935.5

$ws.Range("H65").Value = 3387.1428
$ws.Range("I65").Value = 2986.6667
$ws.Range("J65").Value = 3687.5
$ws.Range("K65").Value = 14933.3335
$ws.Range("L65").Value = 18437.5
$ws.Range("M65").Value = -11813.3335
$ws.Range("N65").Value = -24677.5

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H96").Value = 2319.5715
$ws.Range("I96").Value = 2084
$ws.Range("K96").Value = 2084
$ws.Range("M96").Value = -711

$ws.Range("H126").Value = 12566.353
$ws.Range("I126").Value = 8172.4
$ws.Range("J126").Value = 18843.428
$ws.Range("K126").Value = 24517.2
$ws.Range("L126").Value = 56530.284
$ws.Range("M126").Value = -22047.2
$ws.Range("N126").Value = -61470.284

$ws.Range("H136").Value = 7847.3945
$ws.Range("I136").Value = 1617.4348
$ws.Range("J136").Value = 17400
$ws.Range("K136").Value = 4852.3044
$ws.Range("L136").Value = 52200
$ws.Range("M136").Value = -2302.3044
$ws.Range("N136").Value = -57300

